# Re-generate the "Review date" list/index: every review date whose day-of-month
# is the 15th moves to the 18th (dates on other days-of-month are left as-is).
#
# The dates are stored as plain text (shared strings), not real Excel date
# serials, so a naive `.Value = "yyyy-mm-18"` assignment would get silently
# auto-parsed into a date serial by Excel's "looks like a date" type
# inference - which would also mint a brand new number-format style.
# To keep the cell's literal text type (t="s") and its original (default)
# style untouched, we:
#   1. write the new date as a quoted string FORMULA (a formula string
#      literal is never date-sniffed), then
#   2. Copy / PasteSpecial(xlPasteValues) the cell onto itself, which
#      collapses the formula down to its plain text result in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2

    if ($current -ne $null -and $current.Length -eq 10 -and $current.Substring(8,2) -eq "15") {
        $newValue = $current.Substring(0,8) + "18"

        $cell.Formula = '="' + $newValue + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = $false
